$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

$ws.Range("D2").Value = "61.756.64"
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("D3").Value = "3.385.64"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("E4").Value = "  -0.12%  "
Set-TextValue $ws.Range("D5") "574.83"
$ws.Range("E5").Value = "  +0.26%  "
Set-TextValue $ws.Range("D6") "138.57"
$ws.Range("E6").Value = "  +1.10%  "
Set-TextValue $ws.Range("D7") "1.00"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").Value = "3.383.80"
$ws.Range("E8").Value = "  +0.64%  "
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("E10").Value = "  -1.26%  "
$ws.Range("E11").Value = "  +0.98%  "
Set-TextValue $ws.Range("D12") "0.394"
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("D13").Value = "3.957.40"
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("E14").Value = "  +2.34%  "
$ws.Range("E15").Value = "  +0.27%  "
$ws.Range("E16").Value = "  +3.27%  "
$ws.Range("D17").Value = "3.381.15"
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("D18").Value = "61.834.91"
$ws.Range("E18").Value = "  +0.76%  "
Set-TextValue $ws.Range("D19") "5.93"
$ws.Range("E19").Value = "  +0.64%  "
Set-TextValue $ws.Range("D20") "14.02"
$ws.Range("E20").Value = "  +1.13%  "
Set-TextValue $ws.Range("D21") "9.43"
$ws.Range("E21").Value = "  +0.86%  "
Set-TextValue $ws.Range("D22") "378.98"
$ws.Range("E22").Value = "  -0.19%  "
Set-TextValue $ws.Range("D23") "0.558"
$ws.Range("E23").Value = "  -1.56%  "
$ws.Range("D24").Value = "3.515.33"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("E25").Value = "  -0.19%  "
Set-TextValue $ws.Range("D26") "0.0000126"
$ws.Range("E26").Value = "  +5.21%  "
Set-TextValue $ws.Range("D27") "71.57"
$ws.Range("E27").Value = "  +1.35%  "
Set-TextValue $ws.Range("D28") "1.82"
$ws.Range("E28").Value = "  +11.42%  "
Set-TextValue $ws.Range("D29") "7.67"
$ws.Range("E29").Value = "  -1.30%  "
Set-TextValue $ws.Range("D30") "0.997"
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("E31").Value = "  +4.93%  "
Set-TextValue $ws.Range("D32") "8.33"
$ws.Range("E32").Value = "  +1.41%  "
$ws.Range("E33").Value = "  +1.86%  "
$ws.Range("E34").Value = "  +0.07%  "
Set-TextValue $ws.Range("D35") "23.69"
$ws.Range("E35").Value = "  +1.27%  "
$ws.Range("E36").Value = "  -3.88%  "
$ws.Range("E37").Value = "  -1.98%  "
Set-TextValue $ws.Range("D38") "1.55"
$ws.Range("E38").Value = "  +0.31%  "
Set-TextValue $ws.Range("D39") "164.87"
$ws.Range("E39").Value = "  +2.66%  "
$ws.Range("E40").Value = "  -1.39%  "
$ws.Range("E41").Value = "  +1.50%  "
$ws.Range("E42").Value = "  -0.13%  "
Set-TextValue $ws.Range("D43") "0.776"
$ws.Range("E43").Value = "  +1.83%  "
$ws.Range("E44").Value = "  +1.62%  "
Set-TextValue $ws.Range("D45") "41.66"
$ws.Range("E45").Value = "  +0.45%  "
$ws.Range("E46").Value = "  -0.16%  "
Set-TextValue $ws.Range("D47") "24.47"
$ws.Range("E47").Value = "  +6.17%  "
$ws.Range("E48").Value = "  -0.93%  "
Set-TextValue $ws.Range("D49") "23.14"
$ws.Range("E49").Value = "  +1.76%  "
$ws.Range("D50").Value = "2.370.97"
$ws.Range("E50").Value = "  +1.88%  "
$ws.Range("E51").Value = "  +0.92%  "
